# Generate Report for Archive
# - Status text "Ready for handoff" -> "In Translation" everywhere it appears
#   (Overview!E2:F2/E3:F3, zh-cn!C2:C3, de-de!C2:C3).
# - The corresponding status columns narrow (they were auto-fit to the
#   longer "Ready for handoff" text; re-fit them to the shorter text).

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"

# --- Overview sheet: columns E ("zh-cn") and F ("de-de") hold status text ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# --- zh-cn sheet: column C ("Status") ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

# --- de-de sheet: column C ("Status") ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# --- Narrow the now-shorter status columns (previously auto-fit to
#     "Ready for handoff", now re-fit to "In Translation") ---
$newColumnWidth = 12.576851254417766

$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth
